$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Pre-condition text: "registado" (registered) -> "autenticado" (authenticated)
# ------------------------------------------------------------------
$ws.Range("C4").Value = "Estar autenticado no sistema"

# ------------------------------------------------------------------
# 2. Insert a brand-new scenario step "6. Regista opção" right after the
#    existing "5. Seleciona modelo..." row. This pushes every row below
#    (old rows 12-22) down by one (new rows 13-23), and Excel automatically
#    grows the merged cells that span across the insertion point
#    (B6:B16 -> B6:B17, B17:B19 -> B18:B20, B20:B22 -> B21:B23).
# ------------------------------------------------------------------
$ws.Rows("12:12").Insert()

# Copy the formatting of the row immediately above onto the freshly
# inserted (blank/default-styled) row so it matches the rest of the table.
$ws.Range("B11:D11").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New step text goes in the "System response" column.
$ws.Range("D12").Value = "6. Regista opção"

# ------------------------------------------------------------------
# 3. Renumber every step that came after the newly inserted one (old step
#    "6." becomes "7.", old "7." becomes "8.", old "8." becomes "9.", and
#    the two "Alternativa" sub-flows that used to hang off step 7 now hang
#    off step 8). Also fix the typo "otima" -> "ótima" along the way, and
#    update the "(passo 7)" reference to "(passo 6)" for the second
#    alternative, matching the source content.
# ------------------------------------------------------------------
$ws.Range("D13").Value = "7. Pergunta qual a forma de personalizar o carro"
$ws.Range("C14").Value = "8. Escolhe configuração ótima"
$ws.Range("D15").Value = "9.<<include>> Escolher configuração ótima"

$ws.Range("B18").Value = "Alternativa 1 [Escolher Pacote] (passo 8)"
$ws.Range("C18").Value = "8.1 Escolhe Pacote"
$ws.Range("D19").Value = "8.2 <<include>> Escolhe Pacote"

$ws.Range("B21").Value = "Alternativa 2 [Escolher Especificações] (passo 6)"
$ws.Range("C21").Value = "8.1 Escolhe Especificações"
$ws.Range("D22").Value = "8.2 <<include>> Escolher Especificações"

# ------------------------------------------------------------------
# 4. Match the refreshed view state saved with the workbook.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("D22").Select()
